$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Total_Citations_(2023" row (row 30): update the summary-stat text from
# "4,817 (range: 11 to 102,352)" to "3,752 (range: 0 to 102,352)".
$ws.Range("B30").Value = "3,752 (range: 0 to 102,352)"

# Row 31 (the "    Unknown" row with counts 86 / 0) is being removed entirely.
# Before deleting it, copy its formatting (bottom border, etc.) onto row 30,
# since row 30 becomes the new last row of this block once row 31 is gone.
$ws.Range("A31:E31").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)

# Now delete row 31 outright; rows below (Mean / Wilcoxon / False discovery)
# shift up to become rows 31-33, and their merged ranges shift with them.
$ws.Rows("31").Delete()
